$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (match the header style used by the other header cells:
# bold font, centered/top aligned, thin border all around)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# New data values for columns I (I0) and J (IF)
$iValues = @(9, 4, 8, 10, 1, 9, 9, 7, 8, 8, 10, 7, 6, 8, 7, 9, 7, 5, 6, 5)
$jValues = @(9, 5, 8, 10, 1, 9, 9, 9, 8, 9, 10, 7, 7, 8, 7, 9, 7, 6, 6, 5)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
